$wb = $excel.ActiveWorkbook

# Both the "展览" (Exhibition) sheet and the "全部类型" (All Types) sheet
# carry the same event rows; update the "想去人数" (F column) counts on each.
$ws1 = $wb.Worksheets.Item(1)
$ws4 = $wb.Worksheets.Item(4)

foreach ($ws in @($ws1, $ws4)) {
    $ws.Range("F6").Value = 4577
    $ws.Range("F8").Value = 378
    $ws.Range("F9").Value = 1332
    $ws.Range("F10").Value = 883
    $ws.Range("F12").Value = 944
    $ws.Range("F14").Value = 528
    $ws.Range("F16").Value = 247
}
